# Auto commit at 2026-01-10  7:32:07.34
# Update the Metrics sheet source values (rows 2-13, column B) and move the
# active selections on the "Metrics" and "today" sheets. All dependent
# formulas (on the "today" sheet, and the TODAY()-1 cell) recalculate
# automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 172827.16
$metrics.Range("B3").Value = 130229.87000000001
$metrics.Range("B4").Value = 42218.350000000006
$metrics.Range("B5").Value = 6992
$metrics.Range("B6").Value = 5808697.8899999997
$metrics.Range("B7").Value = 4900947.5
$metrics.Range("B8").Value = 1706310.17
$metrics.Range("B9").Value = 227269
$metrics.Range("B10").Value = 34274078.879999995
$metrics.Range("B11").Value = 32176222.660000004
$metrics.Range("B12").Value = 11988032.209999999
$metrics.Range("B13").Value = 1324899

$metrics.Range("D29").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E9").Select()
